# Apply weekly update: shift existing data rows down by 2 (making room for a
# new week's data at the top of the data block) and append the two rows that
# fall off the bottom as brand-new rows at the end; then populate the freed
# rows 22-23 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Preserve the last two existing data rows (218:219) by copying them to
#    the new rows 220:221 at the end of the sheet. Using Copy(destination)
#    (rather than Copy + PasteSpecial) preserves number formats/styles
#    exactly, including on previously-empty destination cells.
$ws.Range("A218:R219").Copy($ws.Range("A220"))

# 2) Shift the data block rows 22:217 down by two rows, to rows 24:219.
$ws.Range("A22:R217").Copy($ws.Range("A24"))

$excel.CutCopyMode = 0

# 3) Populate the freed rows 22 and 23 with the new week's observations
#    (date 2023-05-04 / serial 45050), Primera and Segunda quality rows.
$ws.Range("D22").Value2 = 45050
$ws.Range("K22").Value2 = 700
$ws.Range("L22").Value2 = 800
$ws.Range("M22").Value2 = 750
$ws.Range("O22").Value2 = "Región de Ñuble"
$ws.Range("P22").Value2 = 750

$ws.Range("D23").Value2 = 45050
$ws.Range("K23").Value2 = 600
$ws.Range("L23").Value2 = 600
$ws.Range("M23").Value2 = 600
$ws.Range("O23").Value2 = "Región de Ñuble"
$ws.Range("P23").Value2 = 600

Write-Output "done"
